$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.288.94"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.605.55"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.90"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.88"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.061.94"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.58"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.276.13"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.606.02"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.34"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.62"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.54"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.535"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.78"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.87"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.41"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("E34").Value = "  +10.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("E36").Value = "  +4.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.64"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "317.43"
$ws.Range("E38").Value = "  +6.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.22"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.88"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.848"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.57"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.94"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.607"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.97"
$ws.Range("E48").Value = "  +5.94%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.97"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.76"
$ws.Range("E51").Value = "  +0.43%  "
